$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 29   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# --- Cells changing type (text <-> number): copy formatting from a stable donor, then set value ---
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Value = 1

$ws.Range("L14").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Value = -100

$ws.Range("I14").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Value = 1

$ws.Range("L14").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value = -100

$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = "0"

$ws.Range("I14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = 1

$ws.Range("L14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = -100

$ws.Range("I14").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").Value = 1

$ws.Range("L14").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").Value = 100

$ws.Range("C14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = "0"

$ws.Range("I14").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").Value = 1

$ws.Range("L14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = -100

$ws.Range("I14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value = 1

$ws.Range("I14").Copy() | Out-Null
$ws.Range("F23").PasteSpecial(-4122) | Out-Null
$ws.Range("F23").Value = 1

$ws.Range("C14").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null
$ws.Range("C26").Value = "0"

$ws.Range("I14").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").Value = 1

$ws.Range("L14").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = -100

# --- Value-only updates (style unchanged) ---
$ws.Range("J14").Value = 13
$ws.Range("K14").Value = -7.692307692307
$ws.Range("N14").Value = -65.714285714285
$ws.Range("J15").Value = 35
$ws.Range("K15").Value = -31.428571428571
$ws.Range("N15").Value = -57.894736842105
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 40
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = 42.857142857142
$ws.Range("I16").Value = 433
$ws.Range("J16").Value = 337
$ws.Range("K16").Value = 28.486646884273
$ws.Range("L16").Value = 52.464788732394
$ws.Range("M16").Value = 12.760416666666
$ws.Range("N16").Value = -60.920577617328
$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 48
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = 33.333333333333
$ws.Range("I17").Value = 603
$ws.Range("J17").Value = 459
$ws.Range("K17").Value = 31.372549019607
$ws.Range("L17").Value = 41.217798594847
$ws.Range("M17").Value = 69.859154929577
$ws.Range("N17").Value = -10.401188707280
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -47.619047619047
$ws.Range("I18").Value = 242
$ws.Range("J18").Value = 212
$ws.Range("K18").Value = 14.150943396226
$ws.Range("L18").Value = 5.217391304347
$ws.Range("M18").Value = 2.978723404255
$ws.Range("N18").Value = -80.530973451327
$ws.Range("C19").Value = 8
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -5.128205128205
$ws.Range("I19").Value = 444
$ws.Range("J19").Value = 441
$ws.Range("K19").Value = 0.680272108843
$ws.Range("L19").Value = 7.766990291262
$ws.Range("M19").Value = 44.625407166123
$ws.Range("N19").Value = 1.138952164009
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 5
$ws.Range("I20").Value = 279
$ws.Range("J20").Value = 190
$ws.Range("K20").Value = 46.842105263157
$ws.Range("L20").Value = 95.104895104895
$ws.Range("M20").Value = 149.107142857143
$ws.Range("N20").Value = -51.730103806228
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -5.263157894736
$ws.Range("F21").Value = 159
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = 8.904109589041
$ws.Range("I21").Value = 2037
$ws.Range("J21").Value = 1687
$ws.Range("K21").Value = 20.746887966805
$ws.Range("L21").Value = 32.444733420026
$ws.Range("M21").Value = 43.754410726887
$ws.Range("N21").Value = -50.713767239293
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = -18.181818181818
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 27
$ws.Range("K23").Value = 50
$ws.Range("L23").Value = -27.027027027027
$ws.Range("M23").Value = 68.75
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 85
$ws.Range("H24").Value = 4.938271604938
$ws.Range("I24").Value = 1164
$ws.Range("J24").Value = 861
$ws.Range("K24").Value = 35.191637630662
$ws.Range("L24").Value = 17.457114026236
$ws.Range("M24").Value = 25.161290322580
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 40
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 54
$ws.Range("H25").Value = 5.555555555555
$ws.Range("I25").Value = 843
$ws.Range("J25").Value = 737
$ws.Range("K25").Value = 14.382632293080
$ws.Range("L25").Value = 9.908735332464
$ws.Range("M25").Value = -9.159482758620
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("J26").Value = 54
$ws.Range("K26").Value = -12.962962962963
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 72
$ws.Range("J27").Value = 57
$ws.Range("K27").Value = 26.315789473684
$ws.Range("L27").Value = 22.033898305084
$ws.Range("C28").Value = 4
$ws.Range("E28").Value = 300
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -14.285714285714
$ws.Range("I28").Value = 54
$ws.Range("J28").Value = 61
$ws.Range("K28").Value = -11.475409836065
$ws.Range("L28").Value = 10.204081632653
$ws.Range("M28").Value = 35
$ws.Range("N28").Value = -59.398496240601
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 45
$ws.Range("J29").Value = 55
$ws.Range("K29").Value = -18.181818181818
$ws.Range("L29").Value = 15.384615384615
$ws.Range("M29").Value = 36.363636363636
$ws.Range("N29").Value = -60.176991150442

Write-Host "Edit complete"